# Auto-generated script to apply scheduled market-data refresh to Anima Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the affected Leve rows
# across all class sheets, matching the upstream scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1044.4482
$ws.Range("I17").Value = 757.8823
$ws.Range("K17").Value = 2273.6469
$ws.Range("M17").Value = -2105.6469
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H64").Value = 3093.3076
$ws.Range("I64").Value = 3028.5715
$ws.Range("J64").Value = 3168.8333
$ws.Range("K64").Value = 3028.5715
$ws.Range("L64").Value = 3168.8333
$ws.Range("M64").Value = -2780.5715
$ws.Range("N64").Value = -3664.8333
$ws.Range("H67").Value = 3093.3076
$ws.Range("I67").Value = 3028.5715
$ws.Range("J67").Value = 3168.8333
$ws.Range("K67").Value = 3028.5715
$ws.Range("L67").Value = 3168.8333
$ws.Range("M67").Value = -2170.5715
$ws.Range("N67").Value = -4884.8333
$ws.Range("H97").Value = 950
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2850
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -3842
$ws.Range("H100").Value = 1900.5555
$ws.Range("I100").Value = 1684.1666
$ws.Range("J100").Value = 2333.3333
$ws.Range("K100").Value = 1684.1666
$ws.Range("L100").Value = 2333.3333
$ws.Range("M100").Value = -1143.1666
$ws.Range("N100").Value = -3415.3333
$ws.Range("H112").Value = 5433.5293
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 5748.125
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 17244.375
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -19460.375
$ws.Range("H113").Value = 2874.9167
$ws.Range("I113").Value = 2897.2
$ws.Range("J113").Value = 2859
$ws.Range("K113").Value = 2897.2
$ws.Range("L113").Value = 2859
$ws.Range("M113").Value = 356.8000000000002
$ws.Range("N113").Value = -9367
$ws.Range("H129").Value = 1272.6028
$ws.Range("I129").Value = 429.42856
$ws.Range("J129").Value = 1613.1154
$ws.Range("K129").Value = 1288.28568
$ws.Range("L129").Value = 4839.3462
$ws.Range("M129").Value = 3711.71432
$ws.Range("N129").Value = -14839.3462
$ws.Range("H138").Value = 2417.4246
$ws.Range("I138").Value = 7319.4
$ws.Range("J138").Value = 2056.9854
$ws.Range("K138").Value = 21958.2
$ws.Range("L138").Value = 6170.956200000001
$ws.Range("M138").Value = -16818.2
$ws.Range("N138").Value = -16450.9562

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 74666.664
$ws.Range("J96").Value = 74666.664
$ws.Range("L96").Value = 74666.664
$ws.Range("N96").Value = -80158.664
$ws.Range("H105").Value = 100000
$ws.Range("J105").Value = 100000
$ws.Range("L105").Value = 100000
$ws.Range("N105").Value = -106988

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 148
$ws.Range("I7").Value = 166.3
$ws.Range("J7").Value = 125.125
$ws.Range("K7").Value = 166.3
$ws.Range("L7").Value = 125.125
$ws.Range("M7").Value = -53.30000000000001
$ws.Range("N7").Value = -351.125
$ws.Range("H22").Value = 100000510
$ws.Range("I22").Value = 353.33334
$ws.Range("J22").Value = 250000750
$ws.Range("K22").Value = 353.33334
$ws.Range("L22").Value = 250000750
$ws.Range("M22").Value = -3.333340000000021
$ws.Range("N22").Value = -250001450
$ws.Range("H62").Value = 2829.8572
$ws.Range("I62").Value = 2758.1667
$ws.Range("J62").Value = 3260
$ws.Range("K62").Value = 2758.1667
$ws.Range("L62").Value = 3260
$ws.Range("M62").Value = -2134.1667
$ws.Range("N62").Value = -4508
$ws.Range("H65").Value = 2829.8572
$ws.Range("I65").Value = 2758.1667
$ws.Range("J65").Value = 3260
$ws.Range("K65").Value = 13790.8335
$ws.Range("L65").Value = 16300
$ws.Range("M65").Value = -10670.8335
$ws.Range("N65").Value = -22540
$ws.Range("H96").Value = 40965.5
$ws.Range("J96").Value = 40965.5
$ws.Range("L96").Value = 40965.5
$ws.Range("N96").Value = -46457.5
$ws.Range("H129").Value = 46999.332
$ws.Range("J129").Value = 46999.332
$ws.Range("L129").Value = 46999.332
$ws.Range("N129").Value = -56999.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1712.8462
$ws.Range("I5").Value = 526.1111
$ws.Range("J5").Value = 2341.1177
$ws.Range("K5").Value = 1578.3333
$ws.Range("L5").Value = 7023.353099999999
$ws.Range("M5").Value = -1466.3333
$ws.Range("N5").Value = -7247.353099999999
$ws.Range("H31").Value = 2150
$ws.Range("I31").Value = 3000
$ws.Range("K31").Value = 9000
$ws.Range("M31").Value = -8712
$ws.Range("H34").Value = 10638797
$ws.Range("I34").Value = 273.6
$ws.Range("J34").Value = 11905288
$ws.Range("K34").Value = 820.8000000000001
$ws.Range("L34").Value = 35715864
$ws.Range("M34").Value = -736.8000000000001
$ws.Range("N34").Value = -35716032
$ws.Range("H120").Value = 11366.363
$ws.Range("J120").Value = 12125
$ws.Range("L120").Value = 36375
$ws.Range("N120").Value = -46051
$ws.Range("H122").Value = 3669
$ws.Range("I122").Value = 468.6842
$ws.Range("J122").Value = 7722.7334
$ws.Range("K122").Value = 4218.1578
$ws.Range("L122").Value = 69504.60060000001
$ws.Range("M122").Value = -1768.1578
$ws.Range("N122").Value = -74404.60060000001
$ws.Range("H131").Value = 3175.0925
$ws.Range("J131").Value = 3576.383
$ws.Range("L131").Value = 10729.149
$ws.Range("N131").Value = -20809.149
$ws.Range("H135").Value = 1712.8462
$ws.Range("I135").Value = 526.1111
$ws.Range("J135").Value = 2341.1177
$ws.Range("K135").Value = 4734.9999
$ws.Range("L135").Value = 21070.0593
$ws.Range("M135").Value = -2199.9999
$ws.Range("N135").Value = -26140.0593
$ws.Range("H137").Value = 6416889.5
$ws.Range("I137").Value = 13898444
$ws.Range("J137").Value = 4128.5713
$ws.Range("K137").Value = 41695332
$ws.Range("L137").Value = 12385.7139
$ws.Range("M137").Value = -41690232
$ws.Range("N137").Value = -22585.7139
$ws.Range("H140").Value = 1385.2106
$ws.Range("I140").Value = 1029.9286
$ws.Range("J140").Value = 2380
$ws.Range("K140").Value = 3089.7858
$ws.Range("L140").Value = 7140
$ws.Range("M140").Value = 2090.2142
$ws.Range("N140").Value = -17500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1091043.2
$ws.Range("I80").Value = 1802678
$ws.Range("J80").Value = 201499.75
$ws.Range("K80").Value = 1802678
$ws.Range("L80").Value = 201499.75
$ws.Range("M80").Value = -1801680
$ws.Range("N80").Value = -203495.75
$ws.Range("H83").Value = 1091043.2
$ws.Range("I83").Value = 1802678
$ws.Range("J83").Value = 201499.75
$ws.Range("K83").Value = 9013390
$ws.Range("L83").Value = 1007498.75
$ws.Range("M83").Value = -9008398
$ws.Range("N83").Value = -1017482.75
$ws.Range("H98").Value = 16540
$ws.Range("J98").Value = 16540
$ws.Range("L98").Value = 16540
$ws.Range("N98").Value = -22530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 127998.625
$ws.Range("I40").Value = 202396.8
$ws.Range("J40").Value = 4001.6667
$ws.Range("K40").Value = 202396.8
$ws.Range("L40").Value = 4001.6667
$ws.Range("M40").Value = -202260.8
$ws.Range("N40").Value = -4273.6667
$ws.Range("H55").Value = 384.68
$ws.Range("I55").Value = 254.66667
$ws.Range("J55").Value = 504.69232
$ws.Range("K55").Value = 254.66667
$ws.Range("L55").Value = 504.69232
$ws.Range("M55").Value = -81.66667000000001
$ws.Range("N55").Value = -850.69232
$ws.Range("H68").Value = 1518.5714
$ws.Range("I68").Value = 1506
$ws.Range("J68").Value = 1550
$ws.Range("K68").Value = 1506
$ws.Range("L68").Value = 1550
$ws.Range("M68").Value = -757
$ws.Range("N68").Value = -3048
$ws.Range("H71").Value = 1518.5714
$ws.Range("I71").Value = 1506
$ws.Range("J71").Value = 1550
$ws.Range("K71").Value = 7530
$ws.Range("L71").Value = 7750
$ws.Range("M71").Value = -3786
$ws.Range("N71").Value = -15238
$ws.Range("H101").Value = 13225.857
$ws.Range("J101").Value = 13225.857
$ws.Range("L101").Value = 13225.857
$ws.Range("N101").Value = -19715.857
$ws.Range("H122").Value = 3591.2727
$ws.Range("I122").Value = 2900.8
$ws.Range("J122").Value = 4166.6665
$ws.Range("K122").Value = 8702.400000000001
$ws.Range("L122").Value = 12499.9995
$ws.Range("M122").Value = -6252.400000000001
$ws.Range("N122").Value = -17399.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 69330
$ws.Range("J94").Value = 69330
$ws.Range("L94").Value = 69330
$ws.Range("N94").Value = -71132
$ws.Range("H95").Value = 79050
$ws.Range("J95").Value = 79050
$ws.Range("L95").Value = 79050
$ws.Range("N95").Value = -84542
$ws.Range("H107").Value = 984
$ws.Range("I107").Value = 1054.091
$ws.Range("J107").Value = 829.8
$ws.Range("K107").Value = 3162.273
$ws.Range("L107").Value = 2489.4
$ws.Range("M107").Value = -1242.273
$ws.Range("N107").Value = -6329.4
$ws.Range("H132").Value = 1619.0938
$ws.Range("I132").Value = 1588.2941
$ws.Range("J132").Value = 1739.9231
$ws.Range("K132").Value = 4764.8823
$ws.Range("L132").Value = 5219.7693
$ws.Range("M132").Value = -2234.8823
$ws.Range("N132").Value = -10279.7693
$ws.Range("H138").Value = 46381.145
$ws.Range("J138").Value = 48933.6
$ws.Range("L138").Value = 48933.6
$ws.Range("N138").Value = -59213.6
$ws.Range("H139").Value = 98985
$ws.Range("J139").Value = 98985
$ws.Range("L139").Value = 98985
$ws.Range("N139").Value = -109265
